$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.049038
$ws.Range("H2").Value = 63.147114
$ws.Range("I2").Value = 0.384846371905728
$ws.Range("J2").Value = 0.384846371905728
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 54.43165466666667
$ws.Range("N2").Value = 163.294964
$ws.Range("O2").Value = 0.2228930782800698
$ws.Range("P2").Value = 0.2327227899462091
$ws.Range("Q2").Value = 1145.733967481544
$ws.Range("R2").Value = 10311.6057073339
$ws.Range("S2").Value = 0.08577959249898426
$ws.Range("T2").Value = 0.0895625213705774

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.049038
$ws.Range("H3").Value = 63.147114
$ws.Range("I3").Value = 0.384846371905728
$ws.Range("J3").Value = 0.384846371905728
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 126.7095336666667
$ws.Range("N3").Value = 380.128601
$ws.Range("O3").Value = 0.5188649542136915
$ws.Range("P3").Value = 0.541747194133123
$ws.Range("Q3").Value = 2667.113789111946
$ws.Range("R3").Value = 24004.02410200752
$ws.Range("S3").Value = 0.1996832951381708
$ws.Range("T3").Value = 0.2084894421522405

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.049038
$ws.Range("H4").Value = 63.147114
$ws.Range("I4").Value = 0.384846371905728
$ws.Range("J4").Value = 0.384846371905728
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 17.491284
$ws.Range("N4").Value = 52.473852
$ws.Range("O4").Value = 0.07162534664261168
$ws.Range("P4").Value = 0.07478406521259567
$ws.Range("Q4").Value = 368.174701584792
$ws.Range("R4").Value = 3313.572314263128
$ws.Range("S4").Value = 0.02756475479189922
$ws.Range("T4").Value = 0.02878037617342881

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.049038
$ws.Range("H5").Value = 63.147114
$ws.Range("I5").Value = 0.384846371905728
$ws.Range("J5").Value = 0.384846371905728
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.62863866666667
$ws.Range("N5").Value = 43.88591599999999
$ws.Range("O5").Value = 0.05990305316690945
$ws.Range("P5").Value = 0.06254481192001105
$ws.Range("Q5").Value = 307.918771182936
$ws.Range("R5").Value = 2771.268940646424
$ws.Range("S5").Value = 0.02305347267736103
$ws.Range("T5").Value = 0.02407014394894238

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 21.049038
$ws.Range("H6").Value = 63.147114
$ws.Range("I6").Value = 0.384846371905728
$ws.Range("J6").Value = 0.384846371905728
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 30.9441155
$ws.Range("N6").Value = 61.888231
$ws.Range("O6").Value = 0.1267135676967176
$ws.Range("P6").Value = 0.08820113878806125
$ws.Range("Q6").Value = 651.3438630358891
$ws.Range("R6").Value = 3908.063178215335
$ws.Range("S6").Value = 0.04876525679931261
$ws.Range("T6").Value = 0.03394388826053896

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 33.645613
$ws.Range("H7").Value = 100.936839
$ws.Range("I7").Value = 0.615153628094272
$ws.Range("J7").Value = 0.615153628094272
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 54.43165466666667
$ws.Range("N7").Value = 163.294964
$ws.Range("O7").Value = 0.2228930782800698
$ws.Range("P7").Value = 0.2327227899462091
$ws.Range("Q7").Value = 1831.386387864311
$ws.Range("R7").Value = 16482.4774907788
$ws.Range("S7").Value = 0.1371134857810855
$ws.Range("T7").Value = 0.1431602685756317

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 33.645613
$ws.Range("H8").Value = 100.936839
$ws.Range("I8").Value = 0.615153628094272
$ws.Range("J8").Value = 0.615153628094272
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 126.7095336666667
$ws.Range("N8").Value = 380.128601
$ws.Range("O8").Value = 0.5188649542136915
$ws.Range("P8").Value = 0.541747194133123
$ws.Range("Q8").Value = 4263.219933159138
$ws.Range("R8").Value = 38368.97939843224
$ws.Range("S8").Value = 0.3191816590755207
$ws.Range("T8").Value = 0.3332577519808825

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 33.645613
$ws.Range("H9").Value = 100.936839
$ws.Range("I9").Value = 0.615153628094272
$ws.Range("J9").Value = 0.615153628094272
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.491284
$ws.Range("N9").Value = 52.473852
$ws.Range("O9").Value = 0.07162534664261168
$ws.Range("P9").Value = 0.07478406521259567
$ws.Range("Q9").Value = 588.5049723370921
$ws.Range("R9").Value = 5296.544751033828
$ws.Range("S9").Value = 0.04406059185071245
$ws.Range("T9").Value = 0.04600368903916686

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 33.645613
$ws.Range("H10").Value = 100.936839
$ws.Range("I10").Value = 0.615153628094272
$ws.Range("J10").Value = 0.615153628094272
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.62863866666667
$ws.Range("N10").Value = 43.88591599999999
$ws.Range("O10").Value = 0.05990305316690945
$ws.Range("P10").Value = 0.06254481192001105
$ws.Range("Q10").Value = 492.1895152955027
$ws.Range("R10").Value = 4429.705637659524
$ws.Range("S10").Value = 0.03684958048954842
$ws.Range("T10").Value = 0.03847466797106867

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 33.645613
$ws.Range("H11").Value = 100.936839
$ws.Range("I11").Value = 0.615153628094272
$ws.Range("J11").Value = 0.615153628094272
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 30.9441155
$ws.Range("N11").Value = 61.888231
$ws.Range("O11").Value = 0.1267135676967176
$ws.Range("P11").Value = 0.08820113878806125
$ws.Range("Q11").Value = 1041.133734740302
$ws.Range("R11").Value = 6246.80240844181
$ws.Range("S11").Value = 0.07794831089740496
$ws.Range("T11").Value = 0.0542572505275223
